# Update "人气" (interest/popularity) counts in column F for the
# "展览" and "全部类型" worksheets, as produced by a regenerated data run.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    5  = 3083
    6  = 1939
    7  = 250
    8  = 77
    9  = 2602
    14 = 144
    15 = 135
    16 = 9775
    18 = 19
    20 = 7723
    21 = 12279
    24 = 254
    25 = 388
    27 = 2760
    29 = 221
    30 = 2785
    31 = 1258
    35 = 4580
    36 = 1231
    37 = 45
    39 = 66
    40 = 601
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "全部类型" (all types) - row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    7  = 3083
    9  = 1939
    11 = 250
    12 = 2602
    18 = 135
    19 = 9775
    20 = 19
    22 = 7723
    23 = 12279
    26 = 254
    30 = 2760
    33 = 221
    37 = 4580
    45 = 601
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
